# Refresh & env added
# Updates the "students" virtual-form sheet:
#  - Row 2: Age (F2) and Mobile No. (I2) switch from text to numeric storage.
#  - Row 3: new student record (Sikhism / EE / General).
#  - Row 4: new student record (Jainism / ECE / ST).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: re-type Age & Mobile No. as numbers (value text stays the same) ---
$ws.Cells.Item(2, 6).Value  = 18        # F2 Age
$ws.Cells.Item(2, 9).Value  = 324243    # I2 Mobile No.

# --- Row 3: new record ---
$ws.Cells.Item(3, 1).Value  = 20240108221909   # A3 Unique id
$ws.Cells.Item(3, 2).Value  = "S"              # B3 First Name
$ws.Cells.Item(3, 3).Value  = "S"              # C3 Last Name
$ws.Cells.Item(3, 4).Value  = "Ms."            # D3 Title
$ws.Cells.Item(3, 5).Value  = "S"              # E3 Guardian's Name
$ws.Cells.Item(3, 6).Value  = 21               # F3 Age
$ws.Cells.Item(3, 7).Value  = "'1/8/24"        # G3 Date of Birth (kept as text)
$ws.Cells.Item(3, 8).Value  = "S"              # H3 Address
$ws.Cells.Item(3, 9).Value  = "S"              # I3 Mobile No.
$ws.Cells.Item(3, 10).Value = "S"              # J3 Email id
$ws.Cells.Item(3, 11).Value = "Male"           # K3 Gender
$ws.Cells.Item(3, 12).Value = "Sikhism"        # L3 Religion
$ws.Cells.Item(3, 13).Value = "General"        # M3 Caste
$ws.Cells.Item(3, 14).Value = "EE"             # N3 Department
$ws.Cells.Item(3, 15).Value = "3rd"            # O3 Year
$ws.Cells.Item(3, 16).Value = "7th"            # P3 Semester
$ws.Cells.Item(3, 17).Value = 1                # Q3 Terms Accepted

# --- Row 4: new record ---
$ws.Cells.Item(4, 1).Value  = 20240108222332   # A4 Unique id
$ws.Cells.Item(4, 2).Value  = "Sayantan"       # B4 First Name
$ws.Cells.Item(4, 3).Value  = "C"              # C4 Last Name
$ws.Cells.Item(4, 4).Value  = "Mr."            # D4 Title
$ws.Cells.Item(4, 5).Value  = "A"              # E4 Guardian's Name
$ws.Cells.Item(4, 6).Value  = "'18"            # F4 Age (kept as text)
$ws.Cells.Item(4, 7).Value  = "'1/8/24"        # G4 Date of Birth (kept as text)
$ws.Cells.Item(4, 8).Value  = "A"              # H4 Address
$ws.Cells.Item(4, 9).Value  = "A"              # I4 Mobile No.
$ws.Cells.Item(4, 10).Value = "A"              # J4 Email id
$ws.Cells.Item(4, 11).Value = "Male"           # K4 Gender
$ws.Cells.Item(4, 12).Value = "Jainism"        # L4 Religion
$ws.Cells.Item(4, 13).Value = "ST"             # M4 Caste
$ws.Cells.Item(4, 14).Value = "ECE"            # N4 Department
$ws.Cells.Item(4, 15).Value = "1st"            # O4 Year
$ws.Cells.Item(4, 16).Value = "1st"            # P4 Semester
$ws.Cells.Item(4, 17).Value = 1                # Q4 Terms Accepted
